$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date in C2
$ws.Range("C2").Value = 45186

# Add the friendly display text "A 60923-2021" as the second HYPERLINK() argument
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_TYRESO/artfynd/A 60923-2021.xlsx", "A 60923-2021")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_TYRESO/kartor/A 60923-2021.png", "A 60923-2021")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_TYRESO/klagomål/A 60923-2021.docx", "A 60923-2021")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_TYRESO/klagomålsmail/A 60923-2021.docx", "A 60923-2021")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_TYRESO/tillsyn/A 60923-2021.docx", "A 60923-2021")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_TYRESO/tillsynsmail/A 60923-2021.docx", "A 60923-2021")'
